{"js": "// Progress Report Edits Continued(1)\n// - Apply Times New Roman to every paragraph/run in the doc body.\n// - Re-flow the subtitle author line into multiple runs, tagging the\n//   mis-spelled surname \"Boyareddigari\" (split \"Boyareddy\"/\"gari\") and the\n//   proper nouns \"Danzer\" with <w:proofErr/> spell-check markers, matching\n//   what Word's background spell checker stamps around words it doesn't\n//   recognize.\n// - Replace the terse \"Issues\" paragraph body with the fuller writeup\n//   (interface/database discussion, crediting Dr. Carini, closing line),\n//   again with proofErr markers around \"Carini\" \u2014 while preserving the\n//   trailing _GoBack bookmark.\n//\n// Because Office.js has no \"set rFonts for the whole body\" shortcut that\n// also inserts <w:proofErr/> markers around specific words, we rebuild each\n// paragraph's contents from an OOXML fragment (wrapped in the minimal\n// pkg:package envelope insertOoxml expects) and splice it in with\n// Range.insertOoxml(..., \"Replace\"). This lets us control run boundaries\n// and proofErr placement exactly, while each paragraph object (and its\n// position in the body) is preserved.\n\nfunction pkg(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships></pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" + bodyXml + \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst RFONTS =\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>';\n\n// Paragraph 1 - Title: \"Issues\"\nconst titleP =\n  \"<w:p><w:pPr><w:pStyle w:val=\\\"Title\\\"/><w:rPr>\" + RFONTS + \"</w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t>Issues</w:t></w:r></w:p>\";\n\n// Paragraph 2 - Subtitle: author byline, split into proofErr-tagged runs.\nconst subtitleP =\n  \"<w:p><w:pPr><w:pStyle w:val=\\\"Subtitle\\\"/><w:rPr>\" + RFONTS + \"</w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t xml:space=\\\"preserve\\\">Karthik </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t>Boyareddy</w:t></w:r>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t>gari</w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t xml:space=\\\"preserve\\\"> and Reese </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t>Danzer</w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellEnd\\\"/></w:p>\";\n\n// Paragraph 3 - the body issues paragraph, with expanded text (interface /\n// database discussion + Dr. Carini credit) and the trailing _GoBack bookmark\n// kept intact at the end of the paragraph.\nconst bodyP =\n  \"<w:p><w:pPr><w:rPr>\" + RFONTS + \"</w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr>\" +\n  \"<w:t>At this point in our project we don\\u2019t possess any major, unsolvable issues</w:t></w:r>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">. We realized that the interface was quite easy to make initially, \" +\n  \"but cannot be completed until we know what data the database will supply. This is something \" +\n  \"that we will collaborate with Dr. </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr><w:t>Carini</w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\"> on to populate the database and create functions to synthesize \" +\n  \"intermediary data points for the display. As of now the interface has temporary values for \" +\n  \"demonstration only; these values will likely need to be replaced with helper functions that \" +\n  \"will feed the correct data to the interface for visualization.</w:t></w:r>\" +\n  \"<w:r><w:rPr>\" + RFONTS + \"</w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\"> However, no coding problems have been thus far encountered.</w:t></w:r>\" +\n  \"<w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/></w:p>\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length < 3) {\n  throw new Error(\"Expected at least 3 paragraphs (title, subtitle, issues body).\");\n}\n\nconst titleRange = paragraphs.items[0].getRange(\"Whole\");\nconst subtitleRange = paragraphs.items[1].getRange(\"Whole\");\nconst bodyRange = paragraphs.items[2].getRange(\"Whole\");\n\ntitleRange.insertOoxml(pkg(titleP), \"Replace\");\nsubtitleRange.insertOoxml(pkg(subtitleP), \"Replace\");\nbodyRange.insertOoxml(pkg(bodyP), \"Replace\");\n\nawait context.sync();\n", "ps1": "# Progress Report Edits Continued(1)\n# - Apply Times New Roman to every paragraph/run in the doc body.\n# - Re-flow the subtitle author line into multiple runs, tagging the\n#   mis-spelled surname \"Boyareddigari\" (split \"Boyareddy\"/\"gari\") and the\n#   proper nouns \"Danzer\" with <w:proofErr/> spell-check markers, matching\n#   what Word's background spell checker stamps around words it doesn't\n#   recognize.\n# - Replace the terse \"Issues\" paragraph body with the fuller writeup\n#   (interface/database discussion, crediting Dr. Carini, closing line),\n#   again with proofErr markers around \"Carini\" -- while preserving the\n#   trailing _GoBack bookmark.\n#\n# The COM object model has no single call that both sets rFonts across a\n# range AND injects <w:proofErr/> markers around individual words, so (same\n# as the Office.js version) each paragraph's contents are rebuilt from an\n# OOXML fragment (wrapped in the minimal pkg:package envelope Range.InsertXML\n# expects) and spliced in over the paragraph's own Range. That gives exact\n# control over run boundaries and proofErr placement while leaving the\n# paragraph objects themselves (and their position in the document) alone.\n\n$d = $word.ActiveDocument\n\nfunction New-Pkg([string]$BodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships></pkg:xmlData></pkg:part>' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $BodyXml + '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$RFONTS = '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>'\n\n# Paragraph 1 - Title: \"Issues\"\n$titleP = '<w:p><w:pPr><w:pStyle w:val=\"Title\"/><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Issues</w:t></w:r></w:p>'\n\n# Paragraph 2 - Subtitle: author byline, split into proofErr-tagged runs.\n$subtitleP = '<w:p><w:pPr><w:pStyle w:val=\"Subtitle\"/><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t xml:space=\"preserve\">Karthik </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Boyareddy</w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>gari</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t xml:space=\"preserve\"> and Reese </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Danzer</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/></w:p>'\n\n# Paragraph 3 - the body issues paragraph, with expanded text (interface /\n# database discussion + Dr. Carini credit) and the trailing _GoBack bookmark\n# kept intact at the end of the paragraph.\n$bodyP = '<w:p><w:pPr><w:rPr>' + $RFONTS + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr>' +\n    '<w:t>At this point in our project we don' + [char]0x2019 + 't possess any major, unsolvable issues</w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr>' +\n    '<w:t xml:space=\"preserve\">. We realized that the interface was quite easy to make initially, but cannot be completed until we know what data the database will supply. This is something that we will collaborate with Dr. </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Carini</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr>' +\n    '<w:t xml:space=\"preserve\"> on to populate the database and create functions to synthesize intermediary data points for the display. As of now the interface has temporary values for demonstration only; these values will likely need to be replaced with helper functions that will feed the correct data to the interface for visualization.</w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr>' +\n    '<w:t xml:space=\"preserve\"> However, no coding problems have been thus far encountered.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n\nif ($d.Paragraphs.Count -lt 3) {\n    throw \"Expected at least 3 paragraphs (title, subtitle, issues body).\"\n}\n\n$titleRange = $d.Paragraphs(1).Range()\n$titleRange.InsertXML((New-Pkg $titleP))\n\n$subtitleRange = $d.Paragraphs(2).Range()\n$subtitleRange.InsertXML((New-Pkg $subtitleP))\n\n$bodyRange = $d.Paragraphs(3).Range()\n$bodyRange.InsertXML((New-Pkg $bodyP))\n"}
